$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.594.24'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.644.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.90'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.625'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.127'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +3.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.395'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.77'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -2.42%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.65'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000195'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.115.60'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.450.50'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.624.73'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.56'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.54%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.42'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.99%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '348.02'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.79'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000112'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.60'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.64'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.18%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.58'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.164'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '534.48'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.82'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.33%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -2.27%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.75'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.38'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -2.25%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '20.27'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.91'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.04%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '153.18'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -3.95%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '159.83'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -3.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.05'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.85%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.29'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0602'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.35%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.48'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -2.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.633'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -2.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0254'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0990'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0252'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.38%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.58'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.94%  '
